$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BE holds the "29-ago" series, appended right after the existing
# "28-ago" column (BD). Mirror BD's formatting (number format + alignment)
# onto the new column before writing values so the new cells land on the
# same style as the rest of the date columns.
$ws.Range("BE1").NumberFormat = $ws.Range("BD1").NumberFormat
$ws.Range("BE2:BE11").NumberFormat = $ws.Range("BD2").NumberFormat
$ws.Range("BE2:BE11").HorizontalAlignment = $ws.Range("BD2").HorizontalAlignment

$ws.Range("BE1").Value = "29-ago"

$ws.Range("BE2").Value = 13
$ws.Range("BE3").Value = 8
$ws.Range("BE4").Value = 9
$ws.Range("BE5").Value = 10
$ws.Range("BE6").Value = 10
$ws.Range("BE7").Value = 15
$ws.Range("BE8").Value = 10
$ws.Range("BE9").Value = 15
$ws.Range("BE10").Value = 19
$ws.Range("BE11").Value = 10

# Match the author's cursor position left after the edit.
$ws.Range("BM10").Select() | Out-Null
